$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 41 --------------------------------------------------------------
# Copy formatting down from the last existing data row (40) so the new
# row picks up the same cell styles (date / time / wrap-text styles).
$ws.Range("E40:M40").Copy()
$ws.Range("E41:M41").PasteSpecial(-4122)

$ws.Range("E41").Value = 44273
$ws.Range("F41").Value = 0.4861111111111111
$ws.Range("G41").Value = 0.51041666666666663
$ws.Range("H41").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"
$ws.Range("I41").Value = "Documentation"
$ws.Range("J41").Value = "Renommer les isseues"
$ws.Range("K41").Value = "CPNV"
$ws.Range("L41").Value = "Rennommer toutes les issues avec la methode smart"
$ws.Range("M41").Value = ""
$ws.Rows.Item(41).RowHeight = 43.2

# --- Row 42 --------------------------------------------------------------
$ws.Range("E40:M40").Copy()
$ws.Range("E42:M42").PasteSpecial(-4122)

$ws.Range("E42").Value = 44274
$ws.Range("F42").Value = 0.58333333333333337
$ws.Range("G42").Value = 0.61458333333333337
$ws.Range("H42").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"
$ws.Range("I42").Value = "Développement"
$ws.Range("J42").Value = "inscription du score"
$ws.Range("K42").Value = "CPNV"
$ws.Range("L42").Value = "Inscription du score d'un un fichier externe"
$ws.Range("M42").Value = "https://www.tutorialspoint.com/c_standard_library/c_function_sp`nrintf.htm`nhttps://codeforwin.org/2018/02/c-program-append-data-file.html"
$ws.Rows.Item(42).RowHeight = 158.4

# --- Grow the table / autofilter to cover the two new rows ---------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("E5:M42"))

# --- Update selection to mirror the author's final cursor position -------
$ws.Range("M43").Select()

Write-Host "done"
